# Weekly update: insert the newest week of "Pepino ensalada" price data
# (Agrícola del Norte S.A. de Arica) at the top of the existing series.
# This pushes the prior rows down by two (rows 96-191 -> 98-193) and the
# two oldest rows that fell off the bottom of the visible range reappear
# as the new rows 192-193 (they were already the last two rows, so the
# insert just relocates them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 96, pushing everything from row 96
# down to row 98 (and so on, through the end of the sheet).
$ws.Rows.Item(96).Resize(2).Insert()

# --- New row 96: Primera ---
$ws.Cells.Item(96,1).Value = 1
$ws.Cells.Item(96,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(96,3).Value = "Arica y Parinacota"
$ws.Cells.Item(96,4).Value = 44484
$ws.Cells.Item(96,5).Value = 15
$ws.Cells.Item(96,6).Value = 100112043
$ws.Cells.Item(96,7).Value = "Pepino ensalada"
$ws.Cells.Item(96,8).Value = "Sin especificar"
$ws.Cells.Item(96,9).Value = "Primera"
$ws.Cells.Item(96,10).Value = 160
$ws.Cells.Item(96,11).Value = 8500
$ws.Cells.Item(96,12).Value = 9000
$ws.Cells.Item(96,13).Value = 8750
$ws.Cells.Item(96,14).Value = "$/caja 70 unidades"
$ws.Cells.Item(96,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(96,16).Value = 125
$ws.Cells.Item(96,17).Value = 70
$ws.Cells.Item(96,18).Value = "Hortaliza"

# --- New row 97: Segunda ---
$ws.Cells.Item(97,1).Value = 1
$ws.Cells.Item(97,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(97,3).Value = "Arica y Parinacota"
$ws.Cells.Item(97,4).Value = 44484
$ws.Cells.Item(97,5).Value = 15
$ws.Cells.Item(97,6).Value = 100112043
$ws.Cells.Item(97,7).Value = "Pepino ensalada"
$ws.Cells.Item(97,8).Value = "Sin especificar"
$ws.Cells.Item(97,9).Value = "Segunda"
$ws.Cells.Item(97,10).Value = 160
$ws.Cells.Item(97,11).Value = 7000
$ws.Cells.Item(97,12).Value = 7500
$ws.Cells.Item(97,13).Value = 7250
$ws.Cells.Item(97,14).Value = "$/caja 100 unidades"
$ws.Cells.Item(97,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(97,16).Value = 72
$ws.Cells.Item(97,17).Value = 100
$ws.Cells.Item(97,18).Value = "Hortaliza"
